# ZEV biosample sheet correction:
#  - treatment column (J): "EtOH" control samples are actually "mockEstradiol"
#  - replicate column (L): TYE7 rows (2-17) use replicate numbers 5/6 instead of 1/2
#  - unhide the previously-hidden helper columns (D, E, F, H, I, J)
#  - update the active selection left over from editing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the treatment column: every "EtOH" cell becomes "mockEstradiol" ---
for ($row = 2; $row -le 57; $row++) {
    $cell = $ws.Cells.Item($row, 10)  # column J = treatment
    if ($cell.Value2 -eq "EtOH") {
        $cell.Value = "mockEstradiol"
    }
}

# --- Fix the replicate numbers for the TYE7 block (rows 2-17) ---
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 12).Value = 5   # column L = replicate
}
for ($row = 10; $row -le 17; $row++) {
    $ws.Cells.Item($row, 12).Value = 6   # column L = replicate
}

# --- Unhide the helper columns that were hidden for the TYE7-only view ---
$ws.Columns.Item(4).Hidden = $false   # D experimentDesign
$ws.Columns.Item(5).Hidden = $false   # E experimentObservations
$ws.Columns.Item(6).Hidden = $false   # F strain
$ws.Columns.Item(8).Hidden = $false   # H floodmedia
$ws.Columns.Item(9).Hidden = $false   # I inductionDelay
$ws.Columns.Item(10).Hidden = $false  # J treatment

# --- Update the lingering selection to a single cell ---
$ws.Range("M12").Select() | Out-Null
